$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows with refreshed market data.
# Price cells hold numeric-looking text (e.g. "91.785.06", "0.0000255") in the source
# workbook, so force a text number-format before assignment and clear it again
# afterwards so Excel keeps the literal string instead of coercing it to a float.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "91.785.06"
$c.ClearFormats()
$ws.Range("E2").Value = "  +0.58%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.119.18"
$c.ClearFormats()
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  -0.09%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "243.65"
$c.ClearFormats()
$ws.Range("E5").Value = "  +2.93%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "626.36"
$c.ClearFormats()
$ws.Range("E6").Value = "  -2.50%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.14"
$c.ClearFormats()
$ws.Range("E7").Value = "  +6.23%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.373"
$c.ClearFormats()
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  -0.08%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.117.33"
$c.ClearFormats()
$ws.Range("E10").Value = "  -1.07%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.769"
$c.ClearFormats()
$ws.Range("E11").Value = "  +5.85%  "
$ws.Range("E12").Value = "  +3.52%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000255"
$c.ClearFormats()
$ws.Range("E13").Value = "  +1.28%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "35.65"
$c.ClearFormats()
$ws.Range("E14").Value = "  -2.24%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "91.658.52"
$c.ClearFormats()
$ws.Range("E15").Value = "  +0.98%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.49"
$c.ClearFormats()
$ws.Range("E16").Value = "  -2.09%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.700.43"
$c.ClearFormats()
$ws.Range("E17").Value = "  -0.95%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.135.95"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.29%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.76"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.13%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0000222"
$c.ClearFormats()
$ws.Range("E20").Value = "  +2.17%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.72"
$c.ClearFormats()
$ws.Range("E21").Value = "  +1.78%  "
$ws.Range("E22").Value = "  +2.93%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "447.10"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.81%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "9.15"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.77%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "5.88"
$c.ClearFormats()
$ws.Range("E25").Value = "  +1.22%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "91.92"
$c.ClearFormats()
$ws.Range("E26").Value = "  +0.33%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.98"
$c.ClearFormats()
$ws.Range("E27").Value = "  -4.10%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  +0.04%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "9.34"
$c.ClearFormats()
$ws.Range("E32").Value = "  -4.41%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.117"
$c.ClearFormats()
$ws.Range("E33").Value = "  +37.03%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.02"
$c.ClearFormats()
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("E35").Value = "  +11.43%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "26.64"
$c.ClearFormats()
$ws.Range("E36").Value = "  -1.53%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "7.67"
$c.ClearFormats()
$ws.Range("E37").Value = "  +7.37%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.13"
$c.ClearFormats()
$ws.Range("E38").Value = "  +22.46%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "492.16"
$c.ClearFormats()
$ws.Range("E41").Value = "  -4.46%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.29"
$c.ClearFormats()
$ws.Range("E42").Value = "  -1.81%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.425"
$c.ClearFormats()
$ws.Range("E43").Value = "  +0.54%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "22.15"
$c.ClearFormats()
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E48").Value = "  -1.38%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "4.59"
$c.ClearFormats()
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  -0.33%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "44.78"
$c.ClearFormats()
$ws.Range("E51").Value = "  -1.87%  "

# Rows 30/31, 39/40 and 46/47 swapped rank order; refresh coin name, link, price and volume.
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.186"
$c.ClearFormats()
$ws.Range("E30").Value = "  +15.95%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.243"
$c.ClearFormats()
$ws.Range("E31").Value = "  +20.28%  "
$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.92"
$c.ClearFormats()
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.65"
$c.ClearFormats()
$ws.Range("E40").Value = "  -5.33%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "156.68"
$c.ClearFormats()
$ws.Range("E46").Value = "  +3.39%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.92"
$c.ClearFormats()
$ws.Range("E47").Value = "  -1.85%  "
